# Build the "購入品リスト1" sheet from the "template" sheet and fill it in
# with the first few parts pulled from the (external) parts list.
$wb = $excel.ActiveWorkbook
$template = $wb.Worksheets.Item("template")

# Duplicate the template so all formatting / formulas / styles carry over,
# then rename and re-position the copy right after "template".
$template.Copy($null, $template)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "購入品リスト1"

# Populate the SKU + quantity columns read from the parts list; the
# subtotal/total formulas already on the sheet pick these up automatically.
$newSheet.Range("A3").Value = 114659
$newSheet.Range("D3").Value = 1

$newSheet.Range("A4").Value = 129604
$newSheet.Range("D4").Value = 2

$newSheet.Range("A5").Value = 109848
$newSheet.Range("D5").Value = 3

# Make the new sheet the active tab.
$newSheet.Activate()

# The cursor on the template sheet had moved to B11 before the switch.
$template.Range("B11").Select()
$newSheet.Activate()
